$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "Arduino Nano" note (row 3, column E) with the extra
#        sentence about the clone board, and wrap/resize the row.
$newArduinoNote = "Cheaper third-party alternatives are available, if you trust them." + [char]10 + "I personally used this clone board: https://www.amazon.com/dp/B09KGVDXZY"
$ws.Cells.Item(3, 5).Value = $newArduinoNote
$ws.Cells.Item(3, 5).WrapText = $true
$ws.Rows(3).RowHeight = 45

# --- 2. Insert a new row for the "2.54mm Breakaway Header/Pins" component
#        right after the 4-pin waterproof connector (row 16), which pushes
#        the "Optional" section (and everything below it) down by one row.
$ws.Rows(17).Insert()

$ws.Cells.Item(17, 1).Value = "2.54mm Breakaway Header/Pins"
$ws.Cells.Item(17, 2).Value = 1
$ws.Cells.Item(17, 3).Value = "A1, U2, U3"
$ws.Cells.Item(17, 4).Value = "https://www.amazon.com/gp/product/B08DVGCTKT"
$ws.Cells.Item(17, 5).Value = "Used to socket the Arduino and also to mount the current sensors. For current sensors, you will use 4 segments of 3-pin connector. For two of them, remove the middle pin. You'll see where they go."
$ws.Rows(17).RowHeight = 45

# --- 3. Update selection to match the new "Optional" header row location.
$ws.Range("A18:D18").Select()
